$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two "Weekend" rows (B5/C5 and B6/C6) and the work entry row (B7/C7)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "Weekend"

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "Weekend"

$ws.Range("B7").Value = 8
$ws.Range("C7").Value = "Finished comments & summaries, added Doxygen code documentation to Git repository. Emailed Riemer about the direction of the project."

# Widen column C to fit the new, longer description text
$ws.Columns.Item(3).ColumnWidth = 127.8

# Update the active selection to C6, matching the author's last-edited cell
$ws.Range("C6").Select()
